$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: a new price record is prepended above the current first
# detail row (row 281), pushing the existing rows down by one.
$ws.Rows.Item(281).Insert()

$ws.Cells.Item(281, 1).Value = 9
$ws.Cells.Item(281, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(281, 3).Value = "Metropolitana"
$ws.Cells.Item(281, 4).Value = 44776
$ws.Cells.Item(281, 5).Value = 13
$ws.Cells.Item(281, 6).Value = 100112021
$ws.Cells.Item(281, 7).Value = "Ají"
$ws.Cells.Item(281, 8).Value = "Americana (o)"
$ws.Cells.Item(281, 9).Value = "Primera"
$ws.Cells.Item(281, 10).Value = 16
$ws.Cells.Item(281, 11).Value = 52000
$ws.Cells.Item(281, 12).Value = 55000
$ws.Cells.Item(281, 13).Value = 53500
$ws.Cells.Item(281, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(281, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(281, 16).Value = 2140
$ws.Cells.Item(281, 17).Value = 25
$ws.Cells.Item(281, 18).Value = "Hortaliza"
